$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column K ("Statut de la domiciliation"),
# shifting it (and everything after it) one column to the right.
$ws.Columns("K:K").Insert()

# Populate the header of the newly inserted column K on row 2.
$rng = $ws.Range("K2")
$rng.Value2 = "Numéro de distribution spéciale (BP, TSA, etc)"

# Match the bold header formatting used by the rest of row 2, with the
# cell's own vertically centered, wrapped text alignment.
$rng.Font.Name = "Calibri"
$rng.Font.Size = 12
$rng.Font.Bold = $true
$rng.VerticalAlignment = -4108
$rng.WrapText = $true

# Give the new column a sensible custom (non best-fit) width, matching the
# sizing used for the rest of the sheet's headers.
$ws.Columns("K:K").ColumnWidth = 21

# Reflect the new edit location in the frozen pane's active selection.
$ws.Range("K5").Select() | Out-Null
